$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 72; this shifts the existing rows 72-79
# down to 74-81, keeping all their data/formatting intact.
$ws.Range("A72:R73").EntireRow.Insert()

# Row 72: new weekly entry (Primera)
$ws.Cells.Item(72, 1).Value = 9
$ws.Cells.Item(72, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(72, 3).Value = "Metropolitana"
$ws.Cells.Item(72, 4).Value = 44578
$ws.Cells.Item(72, 5).Value = 13
$ws.Cells.Item(72, 6).Value = 100114007
$ws.Cells.Item(72, 7).Value = "Jengibre"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 790
$ws.Cells.Item(72, 11).Value = 11000
$ws.Cells.Item(72, 12).Value = 12000
$ws.Cells.Item(72, 13).Value = 11494
$ws.Cells.Item(72, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(72, 15).Value = "Perú"
$ws.Cells.Item(72, 16).Value = 884
$ws.Cells.Item(72, 17).Value = 13
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# Row 73: new weekly entry (Segunda)
$ws.Cells.Item(73, 1).Value = 9
$ws.Cells.Item(73, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(73, 3).Value = "Metropolitana"
$ws.Cells.Item(73, 4).Value = 44578
$ws.Cells.Item(73, 5).Value = 13
$ws.Cells.Item(73, 6).Value = 100114007
$ws.Cells.Item(73, 7).Value = "Jengibre"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Segunda"
$ws.Cells.Item(73, 10).Value = 340
$ws.Cells.Item(73, 11).Value = 10000
$ws.Cells.Item(73, 12).Value = 10000
$ws.Cells.Item(73, 13).Value = 10000
$ws.Cells.Item(73, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(73, 15).Value = "Perú"
$ws.Cells.Item(73, 16).Value = 769
$ws.Cells.Item(73, 17).Value = 13
$ws.Cells.Item(73, 18).Value = "Hortaliza"
